# Update ACR sheet (Sheet2): rows 1-3 split at AX/AY, rows 4-5 split at AW/AX
$wb = $excel.ActiveWorkbook

$acr = $wb.Worksheets.Item("ACR")
$acr.Range("A1:AX3").Value = 0.23255813953488372
$acr.Range("AY1:CW3").Value = 0.76744186046511631
$acr.Range("A4:AW5").Value = 0.23255813953488372
$acr.Range("AX4:CW5").Value = 0.76744186046511631

# Update FAR sheet (Sheet3): AX4 and AX5 change from 1 to 0
$far = $wb.Worksheets.Item("FAR")
$far.Range("AX4:AX5").Value = 0

# Update FRR sheet (Sheet4): AX4 and AX5 change from 0 to 1
$frr = $wb.Worksheets.Item("FRR")
$frr.Range("AX4:AX5").Value = 1
